# Map032 scene update: insert the English text currently duplicated in
# column C into a new column D, and clear out column C so that each row's
# B (Japanese/base) value is followed by an empty C and the value lands in D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $cValue = $cCell.Value2
    if ($null -ne $cValue -and $cValue -ne "") {
        $ws.Cells.Item($r, 4).Value = $cValue
        $cCell.ClearContents()
    }
}
